$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item(1).Name = "summ23605656"
$wb.Worksheets.Item(2).Name = "summ23703736"
$wb.Worksheets.Item(3).Name = "summ23807297"
$wb.Worksheets.Item(4).Name = "summ23915319"
$wb.Worksheets.Item(5).Name = "summ24074772"
$wb.Worksheets.Item(6).Name = "summ24188981"
$wb.Worksheets.Item(7).Name = "summ24284170"
$wb.Worksheets.Item(8).Name = "summ24377273"
$wb.Worksheets.Item(9).Name = "summ24468357"
$wb.Worksheets.Item(10).Name = "summ24560755"
$wb.Worksheets.Item(11).Name = "summ24654609"
$wb.Worksheets.Item(12).Name = "summ24747535"
$wb.Worksheets.Item(13).Name = "summ24841916"
$wb.Worksheets.Item(14).Name = "summ24936363"
$wb.Worksheets.Item(15).Name = "summ25031067"
$wb.Worksheets.Item(16).Name = "summ25123927"
$wb.Worksheets.Item(17).Name = "summ25217785"
$wb.Worksheets.Item(18).Name = "summ25313436"
$wb.Worksheets.Item(19).Name = "summ25404175"
$wb.Worksheets.Item(20).Name = "summ25537447"
$wb.Worksheets.Item(21).Name = "summ25633230"
$wb.Worksheets.Item(22).Name = "summ25731010"
$wb.Worksheets.Item(23).Name = "summ25822592"
$wb.Worksheets.Item(24).Name = "summ25914309"
$wb.Worksheets.Item(25).Name = "summ26008043"
$wb.Worksheets.Item(26).Name = "summ26103801"
$wb.Worksheets.Item(27).Name = "summ26200435"
$wb.Worksheets.Item(28).Name = "summ26299124"
$wb.Worksheets.Item(29).Name = "summ26394971"
$wb.Worksheets.Item(30).Name = "summ26531603"
$wb.Worksheets.Item(31).Name = "summ26628660"
$wb.Worksheets.Item(32).Name = "summ26725843"
$wb.Worksheets.Item(33).Name = "summ26818413"
$wb.Worksheets.Item(34).Name = "summ26911749"
$wb.Worksheets.Item(35).Name = "summ27007833"
$wb.Worksheets.Item(36).Name = "summ27106304"
$wb.Worksheets.Item(37).Name = "summ27204040"
$wb.Worksheets.Item(38).Name = "summ27297551"
$wb.Worksheets.Item(39).Name = "summ27388407"
$wb.Worksheets.Item(40).Name = "summ27483393"
$wb.Worksheets.Item(41).Name = "summ27581070"
$wb.Worksheets.Item(42).Name = "summ27675114"
$wb.Worksheets.Item(43).Name = "summ27767079"
$wb.Worksheets.Item(44).Name = "summ27857597"
$wb.Worksheets.Item(45).Name = "summ27950482"
$wb.Worksheets.Item(46).Name = "summ28047246"
$wb.Worksheets.Item(47).Name = "summ28205729"
$wb.Worksheets.Item(48).Name = "summ28364739"
$wb.Worksheets.Item(49).Name = "summ28457757"
$wb.Worksheets.Item(50).Name = "summ28553145"
